$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows above the existing row 626, pushing the
# current rows 626-630 down to become rows 631-635 (values/styles move
# with them automatically).
$ws.Range("A626:R630").Insert()

# Populate the 5 newly inserted rows (626-630) with this week's prices.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificacion

$rows = @(
    @{ Row = 626; D = 44939; I = "Extra";   J = 430; K = 2600; L = 2600; M = 2600; N = "`$/unidad"; O = "Región de O'Higgins" ; P = 2600 },
    @{ Row = 627; D = 44939; I = "Primera"; J = 520; K = 2200; L = 2200; M = 2200; N = "`$/unidad"; O = "Región de O'Higgins" ; P = 2200 },
    @{ Row = 628; D = 44939; I = "Segunda"; J = 340; K = 1600; L = 1600; M = 1600; N = "`$/unidad"; O = "Región de O'Higgins" ; P = 1600 },
    @{ Row = 629; D = 44939; I = "Super";   J = 250; K = 3200; L = 3200; M = 3200; N = "`$/unidad"; O = "Región de O'Higgins" ; P = 3200 },
    @{ Row = 630; D = 44939; I = "Tercera"; J = 160; K = 1200; L = 1200; M = 1200; N = "`$/unidad"; O = "Región de O'Higgins" ; P = 1200 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112028
    $ws.Cells.Item($row, 7).Value = "Sandia"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
